# Add the new "Team" sheet as the last sheet in the workbook.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Team"

# Header row (row 1).
$headers = @("player_id", "PA", "AB", "R", "H", "D", "T", "HR", "RBI", "SB", "CS", "W", "K", "HBP", "SH", "SF", "IW", "LOB", "RLSP")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Team batting-stat rows (rows 2-11), one row per player.
$rows = @(
    @("bautj002", 240, 192, 33, 41, 14, 0, 10, 41, 0, 1, 42, 46, 2, 0, 4, 1, 150, 150),
    @("carre001", 186, 160, 32, 46, 6, 1, 5, 15, 5, 2, 18, 42, 2, 5, 1, 0, 111, 111),
    @("donaj001", 351, 290, 66, 86, 20, 1, 21, 58, 5, 1, 53, 58, 6, 1, 1, 5, 201, 201),
    @("encae001", 338, 284, 48, 80, 22, 0, 20, 67, 1, 0, 43, 61, 3, 0, 8, 1, 202, 202),
    @("goinr001", 82, 76, 4, 13, 5, 1, 1, 4, 0, 0, 4, 21, 1, 1, 0, 0, 63, 63),
    @("pillk001", 287, 269, 38, 81, 24, 1, 3, 29, 5, 3, 9, 39, 5, 2, 2, 0, 185, 185),
    @("saunm001", 276, 241, 34, 64, 20, 1, 10, 26, 1, 1, 32, 82, 1, 1, 1, 1, 177, 177),
    @("smoaj001", 187, 164, 24, 37, 4, 0, 10, 21, 0, 0, 21, 62, 2, 0, 0, 1, 126, 126),
    @("travd001", 200, 187, 20, 51, 11, 1, 2, 15, 2, 1, 11, 40, 0, 1, 1, 0, 133, 133),
    @("tulot001", 270, 242, 29, 63, 6, 0, 13, 41, 0, 0, 24, 51, 2, 0, 2, 1, 177, 177)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

# The RBI column (I) carries per-player comments from "Sinto" with a
# corrected figure; the commented cells are bolded.
$commentNotes = @{
    "I2" = "40"
    "I5" = "66"
    "I6" = "5"
    "I8" = "25"
}

foreach ($ref in $commentNotes.Keys) {
    $cell = $ws.Range($ref)
    $cell.Font.Bold = $true
    $cell.AddComment("Sinto:" + "`n" + $commentNotes[$ref])
}

# Restore the selection that was active on the new sheet.
$ws.Range("I5").Select()
